# "Save Setting" feature log entries for Raul's Log (Logs sheet)
# Updates a typo fix on existing rows, a time correction, and appends
# six new log rows for 2016-09-21 (serial 42634).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Existing-row fixes -------------------------------------------------
# NOTE: write order matters here because it controls the order new shared
# strings get appended to the workbook's string table, so we intentionally
# touch C525 (the "2150" time fix) before F520 (the "and" typo fix), then
# fill in the brand new rows 534-539 top-to-bottom, left-to-right.

# Rows 525-529: time corrected from 2200 to 2150 (column C only)
$ws.Cells.Item(525, 3).Value = "2150"
$ws.Cells.Item(526, 3).Value = "2150"
$ws.Cells.Item(527, 3).Value = "2150"
$ws.Cells.Item(528, 3).Value = "2150"
$ws.Cells.Item(529, 3).Value = "2150"

# Rows 520-524: task type corrected from Demo to Setup Skype Kit (column A)
# and the note's typo "abd" -> "and" fixed (column F)
$ws.Cells.Item(520, 1).Value = "Setup Skype Kit"
$ws.Cells.Item(520, 6).Value = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"

$ws.Cells.Item(521, 1).Value = "Setup Skype Kit"
$ws.Cells.Item(521, 6).Value = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"

$ws.Cells.Item(522, 1).Value = "Setup Skype Kit"
$ws.Cells.Item(522, 6).Value = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"

$ws.Cells.Item(523, 1).Value = "Setup Skype Kit"
$ws.Cells.Item(523, 6).Value = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"

$ws.Cells.Item(524, 1).Value = "Setup Skype Kit"
$ws.Cells.Item(524, 6).Value = "Video recording via WinMovie  maker -  web cam and tripod in OSG 1014L"

# --- New rows (534-539), dated 2016-09-21 (serial 42634) ---------------

# Row 534
$ws.Cells.Item(534, 1).Value = "Pickup Mic"
$ws.Cells.Item(534, 2).Value = 42634
$ws.Cells.Item(534, 3).Value = "1900"
$ws.Cells.Item(534, 4).Value = "KT"
$ws.Cells.Item(534, 5).Value = "519"
$ws.Cells.Item(534, 6).Value = "Return mic  to KT 516"

# Row 535
$ws.Cells.Item(535, 1).Value = "Demo"
$ws.Cells.Item(535, 2).Value = 42634
$ws.Cells.Item(535, 3).Value = "1600"
$ws.Cells.Item(535, 4).Value = "ACW"
$ws.Cells.Item(535, 5).Value = "302"
$ws.Cells.Item(535, 6).Value = "demo laptop use"

# Row 536 (no column F)
$ws.Cells.Item(536, 1).Value = "Demo"
$ws.Cells.Item(536, 2).Value = 42634
$ws.Cells.Item(536, 3).Value = "1845"
$ws.Cells.Item(536, 4).Value = "DB"
$ws.Cells.Item(536, 5).Value = "1016"

# Row 537
$ws.Cells.Item(537, 1).Value = "Demo"
$ws.Cells.Item(537, 2).Value = 42634
$ws.Cells.Item(537, 3).Value = "1630"
$ws.Cells.Item(537, 4).Value = "SSB"
$ws.Cells.Item(537, 5).Value = "W141"
$ws.Cells.Item(537, 6).Value = "Using PC, neck mic and posium mic"

# Row 538
$ws.Cells.Item(538, 1).Value = "Operator"
$ws.Cells.Item(538, 2).Value = 42634
$ws.Cells.Item(538, 3).Value = "1700"
$ws.Cells.Item(538, 4).Value = "SSB"
$ws.Cells.Item(538, 5).Value = "W141"
$ws.Cells.Item(538, 6).Value = "Operate event between 17:00-17:45"

# Row 539 (no column F)
$ws.Cells.Item(539, 1).Value = "AV Shutdown"
$ws.Cells.Item(539, 2).Value = 42634
$ws.Cells.Item(539, 3).Value = "2000"
$ws.Cells.Item(539, 4).Value = "SSB"
$ws.Cells.Item(539, 5).Value = "W141"

# --- View state: scroll/selection moved to follow the newly entered data
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 524
$ws.Range("E543").Select()
